$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 4.7308
$ws.Range("B10").Value = 8.589100000000004
$ws.Range("B12").Value = 5.8373
$ws.Range("E13").Value = 11.9988
$ws.Range("B18").Value = 4.887600000000005
$ws.Range("B25").Value = 5.550999999999996
